$d = $word.ActiveDocument

# Locate the "Requisitos" list-bullet paragraph that contains the
# LOB1019 requirement line (currently the last of the three lines) and
# move that line so it becomes the first line of the paragraph, leaving
# the other two lines (LOQ4053, LOB1004) in their original relative order.
$needle = "LOB1019"

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*$needle*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found"
} else {
    $pStart = $target.Range.Start

    # Collect the (relative) offsets of every manual line-break (vertical
    # tab, chr(11)) inside the paragraph; each requirement line in the
    # source document ends with one of these before the next run starts.
    $full = $target.Range.Text
    $breaks = @()
    for ($i = 0; $i -lt $full.Length; $i++) {
        if ([int][char]$full[$i] -eq 11) { $breaks += $i }
    }

    # Find which line (0-based) contains the needle so we know its run
    # boundaries regardless of how many lines/breaks precede it.
    $lineStartRel = 0
    $lineIndex = 0
    for ($b = 0; $b -lt $breaks.Count; $b++) {
        $lineEndRel = $breaks[$b] + 1
        $lineText = $full.Substring($lineStartRel, $lineEndRel - $lineStartRel)
        if ($lineText -like "*$needle*") {
            $lineIndex = $b
            break
        }
        $lineStartRel = $lineEndRel
    }

    $runRelStart = $lineStartRel
    $runRelEnd = $breaks[$lineIndex] + 1

    $runAbsStart = $pStart + $runRelStart
    $runAbsEnd = $pStart + $runRelEnd

    $runRange = $d.Range($runAbsStart, $runAbsEnd)
    $runText = $runRange.Text

    # Insert a copy of that run's text (including its trailing line break)
    # as a brand-new run at the very start of the paragraph -- inserting
    # into a collapsed range at a run boundary creates a distinct w:r
    # rather than merging into the neighbouring run.
    $insertionPoint = $d.Range($pStart, $pStart)
    $insertionPoint.InsertBefore($runText)

    # The original run (and everything after the insertion point) has now
    # shifted forward by the length of the inserted text; delete the now
    #-duplicated original occurrence.
    $shift = $runText.Length
    $origRange = $d.Range($runAbsStart + $shift, $runAbsEnd + $shift)
    $origRange.Delete()
}
